$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The top-most row was an empty "spacer" formatting row (row 1) that pushed
# the real header labels to row 2 and the first data row to row 3. Remove
# it so the header moves to row 1 and the data moves to row 2.
$ws.Rows("1").Delete()

# Re-create the freeze pane at the new header/data boundary (was frozen
# after row 3 / col A, now after row 2 / col A) and refresh the selection.
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("B3").Select()
$excel.ActiveWindow.FreezePanes = $true

# Recolor the header band fill (previously a washed out grey background
# tone, now a deep indigo/blue). Set it (plus the wrap-text that the band
# already had) on A1 and B1 individually, then fan the B1 look out across
# the rest of the header row via a format-only paste so every header cell
# ends up sharing one consistent look.
$ws.Range("A1").Interior.PatternColor = 10040115
$ws.Range("A1").WrapText = $true
$ws.Range("B1").Interior.PatternColor = 10040115
$ws.Range("B1").WrapText = $true
$ws.Range("B1").Copy()
$ws.Range("C1:T1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
